# Apply the "include" flag column (column E) to the climate_watch_emissions sheet.
# Column E indicates whether a SISEPUEDE field is included (1) or purposefully
# excluded (0) to match the Climate Watch method, mirroring the row highlighting
# already present in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Cells.Item(1, 5).Value = "include"
$ws.Cells.Item(1, 5).Font.Bold = $true

# Data rows: 1 = included, 0 = excluded (matches the red/orange-highlighted
# rows already shown in column D)
for ($r = 2; $r -le 10; $r++) { $ws.Cells.Item($r, 5).Value = 1 }
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(13, 5).Value = 0
for ($r = 14; $r -le 21; $r++) { $ws.Cells.Item($r, 5).Value = 1 }
for ($r = 22; $r -le 25; $r++) { $ws.Cells.Item($r, 5).Value = 0 }
for ($r = 26; $r -le 28; $r++) { $ws.Cells.Item($r, 5).Value = 1 }
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(30, 5).Value = 1
$ws.Cells.Item(31, 5).Value = 0
for ($r = 32; $r -le 39; $r++) { $ws.Cells.Item($r, 5).Value = 1 }
for ($r = 40; $r -le 43; $r++) { $ws.Cells.Item($r, 5).Value = 0 }
for ($r = 44; $r -le 46; $r++) { $ws.Cells.Item($r, 5).Value = 1 }
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 5).Value = 1
$ws.Cells.Item(49, 5).Value = 0
for ($r = 50; $r -le 57; $r++) { $ws.Cells.Item($r, 5).Value = 1 }
for ($r = 58; $r -le 61; $r++) { $ws.Cells.Item($r, 5).Value = 0 }
for ($r = 62; $r -le 64; $r++) { $ws.Cells.Item($r, 5).Value = 1 }
$ws.Cells.Item(65, 5).Value = 0
for ($r = 66; $r -le 67; $r++) { $ws.Cells.Item($r, 5).Value = 1 }
for ($r = 68; $r -le 70; $r++) { $ws.Cells.Item($r, 5).Value = 0 }
for ($r = 71; $r -le 160; $r++) { $ws.Cells.Item($r, 5).Value = 1 }

# Restore the cursor/selection to where the author left it
$ws.Range("D11").Select()
